$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'63.846.20"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  -2.21%  '

$ws.Range("D3").Value = "'3.096.77"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  -3.11%  '

$ws.Range("E4").Value = '  -0.17%  '

$ws.Range("D5").Value = "'591.24"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.83%  '

$ws.Range("D6").Value = "'156.86"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +1.62%  '

$ws.Range("E7").Value = '  -0.11%  '

$ws.Range("D8").Value = "'0.539"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -0.48%  '

$ws.Range("D9").Value = "'3.102.94"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -2.92%  '

$ws.Range("D10").Value = "'0.158"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -4.60%  '

$ws.Range("D11").Value = "'5.89"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -2.89%  '

$ws.Range("D12").Value = "'0.451"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -4.54%  '

$ws.Range("D13").Value = "'37.00"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -6.05%  '

$ws.Range("D14").Value = "'0.0000239"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -6.24%  '

$ws.Range("B15").Value = 'TRON'
$ws.Range("C15").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D15").Value = "'0.120"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -1.57%  '

$ws.Range("B16").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C16").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D16").Value = "'3.613.09"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -3.18%  '

$ws.Range("D17").Value = "'7.21"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -3.16%  '

$ws.Range("D18").Value = "'63.834.50"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -1.88%  '

$ws.Range("D19").Value = "'3.106.66"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -2.98%  '

$ws.Range("D20").Value = "'477.27"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -1.36%  '

$ws.Range("D21").Value = "'14.44"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -4.29%  '

$ws.Range("D22").Value = "'0.712"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -7.55%  '

$ws.Range("D23").Value = "'7.56"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -4.33%  '

$ws.Range("D24").Value = "'2.46"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +0.05%  '

$ws.Range("B25").Value = 'Litecoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D25").Value = "'81.29"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -2.92%  '

$ws.Range("B26").Value = 'InternetComputer(DFINITY)'
$ws.Range("C26").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D26").Value = "'12.90"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -6.55%  '

$ws.Range("D27").Value = "'10.58"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +4.83%  '

$ws.Range("E28").Value = '  -0.39%  '

$ws.Range("D29").Value = "'7.51"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -0.21%  '

$ws.Range("D30").Value = "'2.68"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -3.56%  '

$ws.Range("E31").Value = '  +0.11%  '

$ws.Range("E32").Value = '  -4.21%  '

$ws.Range("E33").Value = '  -5.86%  '

$ws.Range("D34").Value = "'27.19"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -4.51%  '

$ws.Range("D35").Value = "'0.0₃0839"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -7.20%  '

$ws.Range("E36").Value = '  -3.07%  '

$ws.Range("D37").Value = "'6.03"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -4.96%  '

$ws.Range("B38").Value = 'Stacks'
$ws.Range("C38").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D38").Value = "'2.26"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -5.39%  '

$ws.Range("B39").Value = 'dogwifhat'
$ws.Range("C39").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D39").Value = "'3.29"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -8.91%  '

$ws.Range("D40").Value = "'50.92"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -1.46%  '

$ws.Range("D41").Value = "'9.21"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -2.80%  '

$ws.Range("D42").Value = "'439.29"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -7.88%  '

$ws.Range("D43").Value = "'0.290"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -4.66%  '

$ws.Range("B44").Value = 'Kaspa'
$ws.Range("C44").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D44").Value = "'0.112"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +0.04%  '

$ws.Range("D45").Value = "'0.0363"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -5.42%  '

$ws.Range("B46").Value = 'Arweave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range("D46").Value = "'40.06"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +2.99%  '

$ws.Range("D47").Value = "'2.826.55"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -4.48%  '

$ws.Range("D48").Value = "'130.84"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -0.88%  '

$ws.Range("D49").Value = "'25.81"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -1.24%  '

$ws.Range("D50").Value = "'0.999"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +0.05%  '

$ws.Range("D51").Value = "'2.25"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -4.45%  '
